# Rename Labrador ("LAB"/"NF"-suffixed) colony names to their corrected
# Newfoundland and Labrador ("NL") equivalents in the ATPU trend data sheet,
# and move the sheet selection/active cell to reflect where the author had
# scrolled/selected next.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Old colony name -> new colony name (column B, "Colony"). Order matches the
# sequence the new names were (re)introduced in the workbook's shared-string
# table so the rebuilt table lines up with the author's saved file.
$colonyRenames = [ordered]@{
    "Bacalhao, LAB"                             = "Bacalhao, NL"
    "Tinker, LAB"                                = "Tinker, NL"
    "Herring Island 1, LAB"                      = "Herring Island 1, NL"
    "Herring Island 2, LAB"                      = "Herring Island 2, NL"
    "Herring Island 3, LAB"                      = "Herring Island 3, NL"
    "North Green, LAB"                           = "North Green, NL"
    "Gannet Clusters 2, LAB"                     = "Gannet Clusters 2, NL"
    "Gannet Clusters 3, LAB"                     = "Gannet Clusters 3, NL"
    "Gannet Clusters 4, LAB"                     = "Gannet Clusters 4, NL"
    "Gannet Clusters 5, LAB"                     = "Gannet Clusters 5, NL"
    "Gannet Clusters 6, LAB"                     = "Gannet Clusters 6, NL"
    "Baccalieu Island, NF"                       = "Baccalieu Island, NL"
    "North Bird Island, NF"                      = "North Bird Island, NL"
    "Pee Pee Island, Witless Bay, NF"            = "Pee Pee Island, NL"
    "Great Island, Witless Bay, NF"              = "Great Island, NL"
    "Gull Island, Witless Bay,NF"                = "Gull Island, NL"
    "Coleman Island, Wadham Islands,NF"          = "Coleman Island, NL"
    "South Penguin Island, Wadham Islands, NF"   = "South Penguin Island, NL"
    "Small Island, Wadham Islands, NF"           = "Small Island, NL"
    "Puffin Islands (Lab), LB"                   = "Puffin Islands, NL"
}

$colonyRange = $ws.Range("B2:B103")
foreach ($oldName in $colonyRenames.Keys) {
    $newName = $colonyRenames[$oldName]
    $colonyRange.Replace($oldName, $newName, 1) | Out-Null
}

# Reflect the author's new scroll position / selection in the saved view:
# pane stays frozen at row 1, window now scrolled so row 77 leads, and the
# active cell/selection moved from F85 to B87.
$win = $excel.ActiveWindow
$win.ScrollRow = 77
$win.ScrollColumn = 1
$ws.Range("B87").Select() | Out-Null
